# Runtime/RuntimeFunctionsAndClasses.xlsx - "Functions" sheet
# WIP edit: toggle which rows the AutoFilter currently hides/shows, and
# extend the AutoFilter criteria (new blank-filter on "Written" column,
# new "Conversion" value on the "Category" column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functions")

# Rows that become hidden (were visible before)
$toHide = @(55, 74, 75, 653, 654, 655, 656, 658, 659, 660, 662, 663, 664, 666, 667, 668, 669, 671, 672, 674, 675, 677, 678, 679, 685, 686, 687, 688, 689, 690, 691, 698, 699)

# Rows that become visible (were hidden before)
$toUnhide = @(87, 98, 99, 104, 169, 171, 172, 216, 220, 313, 314, 385, 426, 427, 731, 733, 775, 776, 777, 778, 792, 914, 915, 916, 917, 920, 926, 927, 928, 929, 931, 983)

foreach ($r in $toHide) {
    $ws.Rows.Item($r).Hidden = $true
}

foreach ($r in $toUnhide) {
    $ws.Rows.Item($r).Hidden = $false
}

# Rebuild the AutoFilter criteria so the resulting <filterColumn> order is
# colId 1 (Written), 4 (Assembly), 5 (Category) - matching column B, E, F.
$rng = $ws.Range("A1:F1199")

# New filter: column B ("Written") - show blanks only
$rng.AutoFilter(2, @(""), 7)

# Existing filter: column E ("Assembly") - Core / VO (unchanged, re-applied
# so it lands after the new column-B filter)
$rng.AutoFilter(5, @("Core", "VO"), 7)

# Updated filter: column F ("Category") - add "Conversion" alongside the
# existing "FixedMemory" value
$rng.AutoFilter(6, @("Conversion", "FixedMemory"), 7)

Write-Output "applied hidden-row toggles and autofilter updates"
